$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; this shifts all existing data rows down by one
# (row 1 -> row 2, row 2 -> row 3, ..., row 63 -> row 64), matching the diff exactly.
$ws.Rows(1).Insert()

# Populate the new header row with the shared-string labels used as column titles.
$ws.Range("A1").Value = "X1"
$ws.Range("B1").Value = "X2"
$ws.Range("C1").Value = "X3"
$ws.Range("D1").Value = "X4"
$ws.Range("E1").Value = "Y"

# Header cells use the same centered style (s="1") as the rest of the data columns.
$ws.Range("A1:E1").HorizontalAlignment = -4108

# Restore the view so it is scrolled back to the top with H8 selected (matches target sheetView).
$ws.Range("H8").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
